# Minor update to mecanum spreadsheet
# - Change Direction(deg) input (D6) from 0 to 90
# - Simplify the normalization formulas in F10/G10/F14/G14 (drop the
#   redundant *SQRT(2) factor) since it is now folded into D18
# - Fold the /SQRT(2) factor into the L1 normalization formula (D18)
# - Update the active selection to I21

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Direction(deg) input value
$ws.Range("D6").Value = 90

# L1 normalization now includes the /SQRT(2) factor
$ws.Range("D18").Formula = "=SQRT(1+TAN(L6))/SQRT(2)"

# Wheel-speed formulas simplified to divide by D18 directly
$ws.Range("F10").Formula = "=E10/D18"
$ws.Range("G10").Formula = "=H10/D18"
$ws.Range("F14").Formula = "=E14/D18"
$ws.Range("G14").Formula = "=H14/D18"

# Match the author's final active-cell selection
$ws.Range("I21").Select()
